$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-match the date format on the destination cell so the upcoming paste
# does not mint a new (redundant) number-format style.
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat

# The current weekly data row (row 16) becomes historical data once this
# week's figures arrive, so push a copy of it down to a new row 17 before
# overwriting row 16 with the new numbers.
$ws.Rows.Item(16).Copy()
$ws.Rows.Item(17).PasteSpecial()

# Update row 16 in place with this week's new values.
$ws.Cells.Item(16, 4).Value = 44714
$ws.Cells.Item(16, 10).Value = 200
$ws.Cells.Item(16, 11).Value = 16000
$ws.Cells.Item(16, 12).Value = 17000
$ws.Cells.Item(16, 13).Value = 16400
$ws.Cells.Item(16, 16).Value = 1093
